# Apply the "Nueva rama con sitio facebook" update:
#  - replace the CP (test case) catalogue with the new Facebook-site list
#  - shrink the used range from A1:L21 down to A1:F12
#  - update row 2's sample data and hyperlink display text
#  - move the active selection to D2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the now-unused columns (G:L) and rows (13:21) first so the
#    remaining writes land on a clean, correctly sized sheet.
$ws.Range("G1:L1").EntireColumn.Delete()
$ws.Range("A13:A21").EntireRow.Delete()

# 2) Header row (unchanged text, just re-asserted).
$ws.Range("A1").Value = "TituloCPs"
$ws.Range("B1").Value = "Dato001"
$ws.Range("C1").Value = "Dato002"
$ws.Range("D1").Value = "Dato003"
$ws.Range("E1").Value = "Dato004"
$ws.Range("F1").Value = "Dato005"

# 3) Row 2 sample data (A2/B2/C2 first, D2 last -- matches shared-string order).
$ws.Range("A2").Value = "CP001_login_fallido"
$ws.Range("B2").Value = "jisola.tsoft@gmail.com"
$ws.Range("C2").Value = 12345678

# 4) New list of test cases (rows 3-12).
$ws.Range("A3").Value = "CP002_login_exitoso"
$ws.Range("A4").Value = "CP003_cerrar_sesion"
$ws.Range("A5").Value = "CP004_modo_oscuro"
$ws.Range("A6").Value = "CP005_buscar_persona"
$ws.Range("A7").Value = "CP006_enviar_solicitud"
$ws.Range("A8").Value = "CP007_cancelar_solicitud"
$ws.Range("A9").Value = "CP008_meGusta_pagina"
$ws.Range("A10").Value = "CP009_crear_publicacion"
$ws.Range("A11").Value = "CP010_crear_historia"
$ws.Range("A12").Value = "CP011_enviar_mensaje"

$ws.Range("D2").Value = "¿Olvidaste tu contraseña?"

# 5) Move the selection like the authored workbook (D2 instead of D3).
$null = $ws.Range("D2").Select()
